$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 176, shifting existing rows 176-222 down to 177-223.
$ws.Rows.Item(176).Insert()

# Populate the new row 176 with the new weekly record.
$ws.Cells.Item(176, 1).Value = 6
$ws.Cells.Item(176, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(176, 3).Value = "Metropolitana"
$ws.Cells.Item(176, 4).Value = 44754
$ws.Cells.Item(176, 5).Value = 13
$ws.Cells.Item(176, 6).Value = 100112029
$ws.Cells.Item(176, 7).Value = "Orégano"
$ws.Cells.Item(176, 8).Value = "Sin especificar"
$ws.Cells.Item(176, 9).Value = "Primera"
$ws.Cells.Item(176, 10).Value = 48
$ws.Cells.Item(176, 11).Value = 17000
$ws.Cells.Item(176, 12).Value = 18000
$ws.Cells.Item(176, 13).Value = 17458
$ws.Cells.Item(176, 14).Value = "$/docena de atados"
$ws.Cells.Item(176, 15).Value = "Región Metropolitana"
$ws.Cells.Item(176, 16).Value = 5819
$ws.Cells.Item(176, 17).Value = 3
$ws.Cells.Item(176, 18).Value = "Hortaliza"
